$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: was "Improve ARIMA existing model" -> becomes the old B3 text
#     ("Dada studying forecasting methods ..."), with the same fill as the
#     "On-Going" header (B1) but WITHOUT the header's bold font.
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Font.Bold = $false
$ws.Range("B2").Value = "Dada studying forecasting methods (https://www.datacamp.com/tutorial/tutorial-time-series-forecasting?irclickid=Qu-WEL35QxyIRzmX30wL5WzCUkD2--zut0y1wg0&irgwc=1&utm_medium=affiliate&utm_source=impact&utm_campaign=1310690#what-is-time-series-forecasting-)"

# --- B3: new task text describing the 2 notebooks, same fill as above plus
#     wrap text, and the row grows taller to fit the new content.
$ws.Range("B1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Font.Bold = $false
$ws.Range("B3").WrapText = $true
$ws.Range("B3").Value = "Create prediction and interpretation of results for Aquifers and Water Spring in 2 separate notebooks as:`n- baseline model (such as e.g. ARIMA) on univariate analysis and get conclusions"

$ws.Rows(3).RowHeight = 72.5

# --- A6 was an empty, styled placeholder cell; clear it out entirely.
$ws.Range("A6").Clear()

$excel.CutCopyMode = $false
